# Generate Report for Handback
#
# The "9c314e72-4676-420c-85ae-d1dd746b7e0c.md" file has been handed back
# from localization, so its status flips from "Ready for handoff" to
# "Handed back: in sync with en-US" on every sheet (Overview, zh-cn, de-de),
# and the "Latest Handback DateTime" column (column G) on each language
# sheet is stamped with the new handback timestamp for that language.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item(2)
$wsZhCn.Range("B3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("G3").Value = "2016-03-10 12:51:57"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item(3)
$wsDeDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("G3").Value = "2016-03-10 12:52:04"
